# Apply the "Que Expenses" workbook edits to the Actual sheet:
#  - "ABS Filament" -> "PETG Filament"
#  - remove the "TPU Filament" line item
#  - add a new "Batch 2" marker (column A) above "Wires"
#  - add a new "Connectors" line item at the end of the list

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Actual")

# Remove the "TPU Filament" row (row 4); everything below shifts up one row
$ws.Rows.Item(4).Delete()

# Rename "ABS Filament" to "PETG Filament"
$ws.Cells.Item(3, 2).Value = "PETG Filament"

# Insert a new row above "Wires" (now row 11) and label it "Batch 2"
$ws.Rows.Item(11).Insert()
$ws.Cells.Item(11, 1).Value = "Batch 2"

# Append a new "Connectors" line item as the final row
$ws.Cells.Item(13, 2).Value = "Connectors"

# Match the saved selection state from the edit
$ws.Range("F13").Select()
